# Kurma_LabExam03Grading.xlsx - grade entry pass
# Fills in the "Points for grading" (column E) awarded scores for the
# "Customer Class" and "Product Class" rubric sections so that each
# criterion receives full marks (matching column D), then moves the
# active selection to E15, mirroring the saved state of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: Customer Class (rows 3-6) ---
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Section 2: Product Class (rows 10-14) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Update the active selection/view to match the saved workbook state.
$ws.Range("E15").Select()
